# Insert a new record row at row 75 (pushing the existing rows 75..101 down to 76..102).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 75:101 down to 76:102 by inserting a new row at 75.
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row 75 with the new record's data.
$ws.Range("A75").Value = 4
$ws.Range("B75").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C75").Value = "Los Lagos"
$ws.Range("D75").Value = 44463
$ws.Range("D75").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E75").Value = 10
$ws.Range("F75").Value = "Fruta"
$ws.Range("G75").Value = 100102
$ws.Range("H75").Value = "Cítricos"
$ws.Range("I75").Value = 100102004
$ws.Range("J75").Value = "Mandarina"
$ws.Range("K75").Value = "Murcott"
$ws.Range("L75").Value = "Primera"
$ws.Range("M75").Value = 600
$ws.Range("N75").Value = 6500
$ws.Range("O75").Value = 6500
$ws.Range("P75").Value = 6500
$ws.Range("Q75").Value = "`$/bandeja 10 kilos"
$ws.Range("R75").Value = "Provincia de Limarí"
$ws.Range("S75").Value = 650
$ws.Range("T75").Value = 10
